$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = @(
  "Họ và tên",
  "Giới tính",
  "Ngày sinh",
  "Email",
  "Số điện thoại",
  "Facebook",
  "Trường đại học",
  "Sinh viên",
  "Chuyên ngành",
  "Thành tích, hoạt động nổi bật",
  "CV1",
  "Có Team",
  "Tên nhóm",
  "Họ và tên 2",
  "Giới tính 2",
  "Ngày sinh 2",
  "Email 2",
  "Số điện thoại 2",
  "Facebook 2",
  "Trường đại học 2",
  "Sinh viên 2",
  "Chuyên ngành 2",
  "Thành tích, hoạt động nổi bật 2",
  "CV2",
  "Họ và tên 3",
  "Giới tính 3",
  "Ngày sinh 3",
  "Email 3",
  "Số điện thoại 3",
  "Facebook 3",
  "Trường đại học 3",
  "Sinh viên 3",
  "Chuyên ngành 3",
  "Thành tích, hoạt động nổi bật 3",
  "CV3"
)

$row2 = @(
  "Sample",
  "Khác",
  "2022-04-23",
  "sample@email.com",
  "00000000000",
  "https://fb.com",
  "Sample University",
  "Năm 1",
  "Sample Major",
  "Sample Achievement",
  "https://gecftu.com/manager/cv/uploads\1650701970888-104052759-SampleCV.pdf",
  "Rồi",
  "Sample Team",
  "Sample 2",
  "Khác",
  "2022-04-23",
  "sample2@email.com",
  "0000000000",
  "https://www.fb.com",
  "Sample University 2",
  "Năm 2",
  "Sample Major 2",
  "Sample achievement",
  "https://gecftu.com/manager/cv/uploads\1650702013769-603311779-SampleCV.pdf",
  "Sample 3",
  "Khác",
  "2022-04-23",
  "sample3@email.com",
  "0000000000",
  "https://fb.com/sample3",
  "Sample University 3",
  "Năm 3",
  "Sample Major 3",
  "Sample achievement",
  "https://gecftu.com/manager/cv/uploads\1650702013781-312483526-SampleCV.pdf"
)

# Expand the used range and force text formatting so that values
# like dates ("2022-04-23") and phone numbers ("00000000000") are
# stored as plain text instead of being auto-converted by Excel.
$range = $ws.Range("A1:AI2")
$range.NumberFormat = "@"

for ($i = 0; $i -lt $row1.Length; $i++) {
  $col = $i + 1
  $ws.Cells.Item(1, $col).Value = $row1[$i]
  $ws.Cells.Item(2, $col).Value = $row2[$i]
}

